$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 755
$ws.Range("I6").Value = 264.2857
$ws.Range("J6").Value = 1442
$ws.Range("K6").Value = 792.8571000000001
$ws.Range("L6").Value = 4326
$ws.Range("M6").Value = -680.8571000000001
$ws.Range("N6").Value = -4550

$ws.Range("H21").Value = 35254.25
$ws.Range("J21").Value = 21000
$ws.Range("L21").Value = 21000
$ws.Range("N21").Value = -21936

$ws.Range("H23").Value = 35254.25
$ws.Range("J23").Value = 21000
$ws.Range("L23").Value = 21000
$ws.Range("N23").Value = -21468

$ws.Range("H135").Value = 726.9091
$ws.Range("I135").Value = 528.7
$ws.Range("K135").Value = 4758.3
$ws.Range("M135").Value = -2223.3

$ws.Range("H138").Value = 4288.8237
$ws.Range("I138").Value = 1366.5834
$ws.Range("J138").Value = 5882.773
$ws.Range("K138").Value = 4099.7502
$ws.Range("L138").Value = 17648.319
$ws.Range("M138").Value = 1040.2498
$ws.Range("N138").Value = -27928.319

$ws.Range("H141").Value = 296218.1
$ws.Range("I141").Value = 1146.742
$ws.Range("J141").Value = 1439619.6
$ws.Range("K141").Value = 3440.226
$ws.Range("L141").Value = 4318858.800000001
$ws.Range("M141").Value = 1739.774
$ws.Range("N141").Value = -4329218.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2874.0588
$ws.Range("I61").Value = 1001.05884
$ws.Range("J61").Value = 4747.0586
$ws.Range("K61").Value = 1001.05884
$ws.Range("L61").Value = 4747.0586
$ws.Range("M61").Value = -789.05884
$ws.Range("N61").Value = -5171.0586

$ws.Range("H136").Value = 2874.0588
$ws.Range("I136").Value = 1001.05884
$ws.Range("J136").Value = 4747.0586
$ws.Range("K136").Value = 3003.17652
$ws.Range("L136").Value = 14241.1758
$ws.Range("M136").Value = -453.17652
$ws.Range("N136").Value = -19341.1758

$ws.Range("H138").Value = 45500
$ws.Range("J138").Value = 45500
$ws.Range("L138").Value = 45500
$ws.Range("N138").Value = -55780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 41500
$ws.Range("J58").Value = 41500
$ws.Range("L58").Value = 41500
$ws.Range("N58").Value = -42088

$ws.Range("H139").Value = 29966.666
$ws.Range("J139").Value = 29966.666
$ws.Range("L139").Value = 29966.666
$ws.Range("N139").Value = -40246.666

$ws.Range("H141").Value = 36415.832
$ws.Range("J141").Value = 29623.75
$ws.Range("L141").Value = 29623.75
$ws.Range("N141").Value = -39983.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2454.3965
$ws.Range("I31").Value = 1561.069
$ws.Range("J31").Value = 3347.724
$ws.Range("K31").Value = 1561.069
$ws.Range("L31").Value = 3347.724
$ws.Range("M31").Value = -1266.069
$ws.Range("N31").Value = -3937.724

$ws.Range("H34").Value = 2454.3965
$ws.Range("I34").Value = 1561.069
$ws.Range("J34").Value = 3347.724
$ws.Range("K34").Value = 1561.069
$ws.Range("L34").Value = 3347.724
$ws.Range("M34").Value = -1359.069
$ws.Range("N34").Value = -3751.724

$ws.Range("H52").Value = 29745
$ws.Range("J52").Value = 29745
$ws.Range("L52").Value = 29745
$ws.Range("N52").Value = -30333

$ws.Range("H58").Value = 10419062
$ws.Range("I58").Value = 1255.125
$ws.Range("J58").Value = 31254676
$ws.Range("K58").Value = 1255.125
$ws.Range("L58").Value = 31254676
$ws.Range("M58").Value = -1052.125
$ws.Range("N58").Value = -31255082

$ws.Range("H99").Value = 3400
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -6996

$ws.Range("H126").Value = 3400
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 2553.775
$ws.Range("I132").Value = 2037.4333
$ws.Range("J132").Value = 4102.8
$ws.Range("K132").Value = 6112.2999
$ws.Range("L132").Value = 12308.4
$ws.Range("M132").Value = -3582.2999
$ws.Range("N132").Value = -17368.4

$ws.Range("H134").Value = 1853.7333
$ws.Range("I134").Value = 880.94116
$ws.Range("J134").Value = 3125.8462
$ws.Range("K134").Value = 2642.82348
$ws.Range("L134").Value = 9377.5386
$ws.Range("M134").Value = -107.82348
$ws.Range("N134").Value = -14447.5386

$ws.Range("H136").Value = 10419062
$ws.Range("I136").Value = 1255.125
$ws.Range("J136").Value = 31254676
$ws.Range("K136").Value = 3765.375
$ws.Range("L136").Value = 93764028
$ws.Range("M136").Value = -1215.375
$ws.Range("N136").Value = -93769128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1249.8334
$ws.Range("I86").Value = 1149.75
$ws.Range("J86").Value = 1450
$ws.Range("K86").Value = 3449.25
$ws.Range("L86").Value = 4350
$ws.Range("M86").Value = -2263.25
$ws.Range("N86").Value = -6722

$ws.Range("H89").Value = 1249.8334
$ws.Range("I89").Value = 1149.75
$ws.Range("J89").Value = 1450
$ws.Range("K89").Value = 10347.75
$ws.Range("L89").Value = 13050
$ws.Range("M89").Value = -4419.75
$ws.Range("N89").Value = -24906

$ws.Range("H107").Value = 1667.3334
$ws.Range("I107").Value = 385
$ws.Range("J107").Value = 1783.909
$ws.Range("K107").Value = 1155
$ws.Range("L107").Value = 5351.727000000001
$ws.Range("M107").Value = 765
$ws.Range("N107").Value = -9191.727000000001

$ws.Range("H122").Value = 1499.45
$ws.Range("J122").Value = 2598.7
$ws.Range("L122").Value = 23388.3
$ws.Range("N122").Value = -28288.3

$ws.Range("H131").Value = 1537.2727
$ws.Range("J131").Value = 1332.1428
$ws.Range("L131").Value = 3996.4284
$ws.Range("N131").Value = -14076.4284

$ws.Range("H133").Value = 3635.6667
$ws.Range("I133").Value = 5792.857
$ws.Range("J133").Value = 2557.0715
$ws.Range("K133").Value = 17378.571
$ws.Range("L133").Value = 7671.2145
$ws.Range("M133").Value = -12318.571
$ws.Range("N133").Value = -17791.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 29546
$ws.Range("J137").Value = 29546
$ws.Range("L137").Value = 29546
$ws.Range("N137").Value = -39746

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2097
$ws.Range("I100").Value = 1456.2
$ws.Range("J100").Value = 2364
$ws.Range("K100").Value = 1456.2
$ws.Range("L100").Value = 2364
$ws.Range("M100").Value = -915.2
$ws.Range("N100").Value = -3446

$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

$ws.Range("H132").Value = 3209.6428
$ws.Range("I132").Value = 1729.1875
$ws.Range("K132").Value = 5187.5625
$ws.Range("M132").Value = -2657.5625

$ws.Range("H136").Value = 1685.6
$ws.Range("I136").Value = 1204.5927
$ws.Range("J136").Value = 2684.6155
$ws.Range("K136").Value = 3613.7781
$ws.Range("L136").Value = 8053.8465
$ws.Range("M136").Value = -1063.7781
$ws.Range("N136").Value = -13153.8465

$ws.Range("H137").Value = 28546.666
$ws.Range("J137").Value = 28546.666
$ws.Range("L137").Value = 28546.666
$ws.Range("N137").Value = -38746.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 15766.667
$ws.Range("I3").Value = 300
$ws.Range("J3").Value = 23500
$ws.Range("K3").Value = 300
$ws.Range("L3").Value = 23500
$ws.Range("M3").Value = -186
$ws.Range("N3").Value = -23728

$ws.Range("H107").Value = 1522.3334
$ws.Range("I107").Value = 462.625
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1387.875
$ws.Range("L107").Value = 30000
$ws.Range("M107").Value = 532.125
$ws.Range("N107").Value = -33840

$ws.Range("H132").Value = 10055.984
$ws.Range("I132").Value = 2073.5186
$ws.Range("J132").Value = 53161.3
$ws.Range("K132").Value = 6220.5558
$ws.Range("L132").Value = 159483.9
$ws.Range("M132").Value = -3690.5558
$ws.Range("N132").Value = -164543.9

$ws.Range("H136").Value = 877.3019
$ws.Range("I136").Value = 529.425
$ws.Range("K136").Value = 1588.275
$ws.Range("M136").Value = 961.7250000000001

$ws.Range("H141").Value = 28369.166
$ws.Range("J141").Value = 28369.166
$ws.Range("L141").Value = 28369.166
$ws.Range("N141").Value = -38729.166

Write-Output "Applied 218 cell updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)"
